$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E18 was stored as text "+917990747606"; it should become a plain number 917990747606
$ws.Range("E18").Value = 917990747606

# Append a new lead as row 19
$ws.Range("A19").Value = "book"
$ws.Range("B19").Value = "Interested"
$ws.Range("C19").Value = "neutral"
$ws.Range("D19").Value = "2025-11-27 17:48:48"

# E19 must stay text (keep the leading "+"), not be auto-coerced to a number
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "+919510038048"
$ws.Range("E19").Style = "Normal"
